$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.656.83"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "2.658.47"
$ws.Range("E3").Value = "  +6.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.176"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +17.01%  "
$ws.Range("D10").Value = "2.657.90"
$ws.Range("E10").Value = "  +6.25%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  +5.56%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000193"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.20%  "
$ws.Range("D15").Value = "3.147.16"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.69%  "
$ws.Range("D17").Value = "72.618.09"
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "2.658.35"
$ws.Range("E18").Value = "  +7.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "384.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.22%  "
$ws.Range("E21").Value = "  +6.06%  "
$ws.Range("E22").Value = "  +5.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +24.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("E25").Value = "  +7.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.29%  "
$ws.Range("D28").Value = "2.797.84"
$ws.Range("E28").Value = "  +6.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "0.0₃0976"
$ws.Range("E30").Value = "  +11.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "545.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.11%  "
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("E39").Value = "  +9.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("E41").Value = "  +10.52%  "
$ws.Range("E42").Value = "  +8.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.29%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  +5.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("E49").Value = "  +7.59%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.17%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0271"
$ws.Range("E51").Value = "  +12.59%  "
